$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Mark the Result column (S) as PASS for the rows that were validated
# (rows 4 and 7 are intentionally left untouched, matching the original data).
foreach ($row in 2, 3, 5, 6, 8) {
    $ws.Range("S$row").Value = "PASS"
}

# Selenium run left the cursor on E13 when the workbook was saved.
$ws.Activate()
$ws.Range("E13").Select()
